# Insert a new weekly price record (row 49) into the Haba sheet, above the
# two previously-last records which get pushed down to rows 50-51 (R50 ->
# R51 in the sheet dimension).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing row 49 (and everything below it) down by one row.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly record.
$ws.Cells.Item(49, 1).Value = 2
$ws.Cells.Item(49, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(49, 3).Value = "Coquimbo"
$ws.Cells.Item(49, 4).Value = 44714
$ws.Cells.Item(49, 5).Value = 4
$ws.Cells.Item(49, 6).Value = 100112026
$ws.Cells.Item(49, 7).Value = "Haba"
$ws.Cells.Item(49, 8).Value = "Sin especificar"
$ws.Cells.Item(49, 9).Value = "Primera"
$ws.Cells.Item(49, 10).Value = 500
$ws.Cells.Item(49, 11).Value = 12000
$ws.Cells.Item(49, 12).Value = 13000
$ws.Cells.Item(49, 13).Value = 12500
$ws.Cells.Item(49, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(49, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(49, 16).Value = 500
$ws.Cells.Item(49, 17).Value = 25
$ws.Cells.Item(49, 18).Value = "Hortaliza"

# Make sure the date column keeps the same date-time number format used by
# the rest of column D (style index 2 in the original workbook).
$ws.Cells.Item(49, 4).NumberFormat = $ws.Cells.Item(50, 4).NumberFormat

Write-Host "Inserted new row 49; rows shifted to 50-51."
